$p = $ppt.ActivePresentation
$sm = $p.SlideMaster
$nm = $p.NotesMaster
Write-Host "before swap"
$tmp = $sm.Theme
$sm.Theme = $nm.Theme
$nm.Theme = $tmp
Write-Host "after swap attempt"
